# Weekly update: insert two new price-report rows for
# "Femacal de La Calera - Ají" right before the existing row 269,
# shifting the remaining rows (old 269-312) down to (271-314).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 269.
$ws.Rows.Item(269).Insert()
$ws.Rows.Item(269).Insert()

# New row 269
$ws.Range("A269").Value = 3
$ws.Range("B269").Value = "Femacal de La Calera"
$ws.Range("C269").Value = "Coquimbo"
$ws.Range("D269").Value = 44522
$ws.Range("E269").Value = 5
$ws.Range("F269").Value = 100112021
$ws.Range("G269").Value = "Ají"
$ws.Range("H269").Value = "Americana (o)"
$ws.Range("I269").Value = "Primera"
$ws.Range("J269").Value = 38
$ws.Range("K269").Value = 26000
$ws.Range("L269").Value = 26000
$ws.Range("M269").Value = 26000
$ws.Range("N269").Value = "$/caja 15 kilos"
$ws.Range("O269").Value = "Limache"
$ws.Range("P269").Value = 1733
$ws.Range("Q269").Value = 15
$ws.Range("R269").Value = "Hortaliza"

# New row 270
$ws.Range("A270").Value = 3
$ws.Range("B270").Value = "Femacal de La Calera"
$ws.Range("C270").Value = "Coquimbo"
$ws.Range("D270").Value = 44522
$ws.Range("E270").Value = 5
$ws.Range("F270").Value = 100112021
$ws.Range("G270").Value = "Ají"
$ws.Range("H270").Value = "Americana (o)"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 71
$ws.Range("K270").Value = 27000
$ws.Range("L270").Value = 28000
$ws.Range("M270").Value = 27493
$ws.Range("N270").Value = "$/caja 15 kilos"
$ws.Range("O270").Value = "Región de Arica y Parinacota"
$ws.Range("P270").Value = 1833
$ws.Range("Q270").Value = 15
$ws.Range("R270").Value = "Hortaliza"

# Apply the same date number-format as the rest of column D so the two
# new date cells match (Excel maps this back onto the existing style).
$ws.Range("D269").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D270").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "Inserted rows 269-270"
